$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Rows whose "handoff generated" timestamps / priority move as part of
# regenerating the handoff report.
$rows = @(7, 9, 10, 12, 13, 14)

foreach ($r in $rows) {
    # Overview!G<r> ("Latest HO Xliff Generate Date") advances from
    # 2016-08-29 14:23:14 to 2016-08-29 14:23:32
    $wsOverview.Range("G$r").Value = "2016-08-29 14:23:32"

    # de-de!H<r> ("Latest Handoff Datetime") shares the same original
    # timestamp text as Overview!G<r>, so it moves together.
    $wsDeDe.Range("H$r").Value = "2016-08-29 14:23:32"

    # zh-cn!H<r> ("Latest Handoff Datetime") advances from
    # 2016-08-29 14:23:06 to 2016-08-29 14:23:27
    $wsZhCn.Range("H$r").Value = "2016-08-29 14:23:27"

    # Priority column (E) newly set to "ht" for both locale sheets.
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"
}
